$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Cfh'
$ws.Cells.Item(2, 3).Value = 'Sell'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 0.6718046666666666
$ws.Cells.Item(2, 8).Value = 2.015414
$ws.Cells.Item(2, 9).Value = 0.006062562609515538
$ws.Cells.Item(2, 10).Value = 0.006345114347523019
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.146719
$ws.Cells.Item(2, 14).Value = 0.440157
$ws.Cells.Item(2, 15).Value = 0.002213357657235064
$ws.Cells.Item(2, 16).Value = 0.002213357657235064
$ws.Cells.Item(2, 17).Value = 0.09856650888866667
$ws.Cells.Item(2, 18).Value = 0.887098579998
$ws.Cells.Item(2, 19).Value = 0.00001341861937423821
$ws.Cells.Item(2, 20).Value = 0.00001404400742712214

# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Cfh'
$ws.Cells.Item(3, 3).Value = 'Sell'
$ws.Cells.Item(3, 4).Value = 'M1'
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 0.6718046666666666
$ws.Cells.Item(3, 8).Value = 2.015414
$ws.Cells.Item(3, 9).Value = 0.006062562609515538
$ws.Cells.Item(3, 10).Value = 0.006345114347523019
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 28.93198366666667
$ws.Cells.Item(3, 14).Value = 86.795951
$ws.Cells.Item(3, 15).Value = 0.4364589970461662
$ws.Cells.Item(3, 16).Value = 0.4364589970461662
$ws.Cells.Item(3, 17).Value = 19.43664164319044
$ws.Cells.Item(3, 18).Value = 174.929774788714
$ws.Cells.Item(3, 19).Value = 0.00264605999607874
$ws.Cells.Item(3, 20).Value = 0.002769382244263136

# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Cfh'
$ws.Cells.Item(4, 3).Value = 'Sell'
$ws.Cells.Item(4, 4).Value = 'M2'
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.6718046666666666
$ws.Cells.Item(4, 8).Value = 2.015414
$ws.Cells.Item(4, 9).Value = 0.006062562609515538
$ws.Cells.Item(4, 10).Value = 0.006345114347523019
$ws.Cells.Item(4, 11).Value = 2.0
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 37.20927366666667
$ws.Cells.Item(4, 14).Value = 111.627821
$ws.Cells.Item(4, 15).Value = 0.5613276452965987
$ws.Cells.Item(4, 16).Value = 0.5613276452965988
$ws.Cells.Item(4, 17).Value = 24.99736369254378
$ws.Cells.Item(4, 18).Value = 224.976273232894
$ws.Cells.Item(4, 19).Value = 0.00340308399406256
$ws.Cells.Item(4, 20).Value = 0.003561688095832761

# Row 5
$ws.Cells.Item(5, 1).Value = 'FAPs'
$ws.Cells.Item(5, 2).Value = 'Cfh'
$ws.Cells.Item(5, 3).Value = 'Sell'
$ws.Cells.Item(5, 4).Value = 'ECs'
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 43.16235166666667
$ws.Cells.Item(5, 8).Value = 129.487055
$ws.Cells.Item(5, 9).Value = 0.3895097374828606
$ws.Cells.Item(5, 10).Value = 0.4076632247761514
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 0.146719
$ws.Cells.Item(5, 14).Value = 0.440157
$ws.Cells.Item(5, 15).Value = 0.002213357657235064
$ws.Cells.Item(5, 16).Value = 0.002213357657235064
$ws.Cells.Item(5, 17).Value = 6.332737074181667
$ws.Cells.Item(5, 18).Value = 56.994633667635
$ws.Cells.Item(5, 19).Value = 0.0008621243600253091
$ws.Cells.Item(5, 20).Value = 0.0009023045201314335

# Row 6
$ws.Cells.Item(6, 1).Value = 'FAPs'
$ws.Cells.Item(6, 2).Value = 'Cfh'
$ws.Cells.Item(6, 3).Value = 'Sell'
$ws.Cells.Item(6, 4).Value = 'M1'
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 43.16235166666667
$ws.Cells.Item(6, 8).Value = 129.487055
$ws.Cells.Item(6, 9).Value = 0.3895097374828606
$ws.Cells.Item(6, 10).Value = 0.4076632247761514
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 28.93198366666667
$ws.Cells.Item(6, 14).Value = 86.795951
$ws.Cells.Item(6, 15).Value = 0.4364589970461662
$ws.Cells.Item(6, 16).Value = 0.4364589970461662
$ws.Cells.Item(6, 17).Value = 1248.772453434923
$ws.Cells.Item(6, 18).Value = 11238.95208091431
$ws.Cells.Item(6, 19).Value = 0.1700050293614848
$ws.Cells.Item(6, 20).Value = 0.1779282822184048

# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Cfh'
$ws.Cells.Item(7, 3).Value = 'Sell'
$ws.Cells.Item(7, 4).Value = 'M2'
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 43.16235166666667
$ws.Cells.Item(7, 8).Value = 129.487055
$ws.Cells.Item(7, 9).Value = 0.3895097374828606
$ws.Cells.Item(7, 10).Value = 0.4076632247761514
$ws.Cells.Item(7, 11).Value = 2.0
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 37.20927366666667
$ws.Cells.Item(7, 14).Value = 111.627821
$ws.Cells.Item(7, 15).Value = 0.5613276452965987
$ws.Cells.Item(7, 16).Value = 0.5613276452965988
$ws.Cells.Item(7, 17).Value = 1606.039755261906
$ws.Cells.Item(7, 18).Value = 14454.35779735715
$ws.Cells.Item(7, 19).Value = 0.2186425837613505
$ws.Cells.Item(7, 20).Value = 0.2288326380376151

# Row 8
$ws.Cells.Item(8, 1).Value = 'M1'
$ws.Cells.Item(8, 2).Value = 'Cfh'
$ws.Cells.Item(8, 3).Value = 'Sell'
$ws.Cells.Item(8, 4).Value = 'ECs'
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 21.686232
$ws.Cells.Item(8, 8).Value = 65.058696
$ws.Cells.Item(8, 9).Value = 0.1957029264426257
$ws.Cells.Item(8, 10).Value = 0.2048238552578966
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 0.146719
$ws.Cells.Item(8, 14).Value = 0.440157
$ws.Cells.Item(8, 15).Value = 0.002213357657235064
$ws.Cells.Item(8, 16).Value = 0.002213357657235064
$ws.Cells.Item(8, 17).Value = 3.181782272808
$ws.Cells.Item(8, 18).Value = 28.636040455272
$ws.Cells.Item(8, 19).Value = 0.0004331605707850961
$ws.Cells.Item(8, 20).Value = 0.0004533484484194719

# Row 9
$ws.Cells.Item(9, 1).Value = 'M1'
$ws.Cells.Item(9, 2).Value = 'Cfh'
$ws.Cells.Item(9, 3).Value = 'Sell'
$ws.Cells.Item(9, 4).Value = 'M1'
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 21.686232
$ws.Cells.Item(9, 8).Value = 65.058696
$ws.Cells.Item(9, 9).Value = 0.1957029264426257
$ws.Cells.Item(9, 10).Value = 0.2048238552578966
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 28.93198366666667
$ws.Cells.Item(9, 14).Value = 86.795951
$ws.Cells.Item(9, 15).Value = 0.4364589970461662
$ws.Cells.Item(9, 16).Value = 0.4364589970461662
$ws.Cells.Item(9, 17).Value = 627.425710015544
$ws.Cells.Item(9, 18).Value = 5646.831390139896
$ws.Cells.Item(9, 19).Value = 0.08541630299414807
$ws.Cells.Item(9, 20).Value = 0.08939721443699067

# Row 10
$ws.Cells.Item(10, 1).Value = 'M1'
$ws.Cells.Item(10, 2).Value = 'Cfh'
$ws.Cells.Item(10, 3).Value = 'Sell'
$ws.Cells.Item(10, 4).Value = 'M2'
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 21.686232
$ws.Cells.Item(10, 8).Value = 65.058696
$ws.Cells.Item(10, 9).Value = 0.1957029264426257
$ws.Cells.Item(10, 10).Value = 0.2048238552578966
$ws.Cells.Item(10, 11).Value = 2.0
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 37.20927366666667
$ws.Cells.Item(10, 14).Value = 111.627821
$ws.Cells.Item(10, 15).Value = 0.5613276452965987
$ws.Cells.Item(10, 16).Value = 0.5613276452965988
$ws.Cells.Item(10, 17).Value = 806.928941286824
$ws.Cells.Item(10, 18).Value = 7262.360471581415
$ws.Cells.Item(10, 19).Value = 0.1098534628776926
$ws.Cells.Item(10, 20).Value = 0.1149732923724865

# Row 11
$ws.Cells.Item(11, 1).Value = 'M2'
$ws.Cells.Item(11, 2).Value = 'Cfh'
$ws.Cells.Item(11, 3).Value = 'Sell'
$ws.Cells.Item(11, 4).Value = 'ECs'
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 30.48803466666666
$ws.Cells.Item(11, 8).Value = 91.46410399999999
$ws.Cells.Item(11, 9).Value = 0.2751329786452017
$ws.Cells.Item(11, 10).Value = 0.2879558237532028
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 0.146719
$ws.Cells.Item(11, 14).Value = 0.440157
$ws.Cells.Item(11, 15).Value = 0.002213357657235064
$ws.Cells.Item(11, 16).Value = 0.002213357657235064
$ws.Cells.Item(11, 17).Value = 4.473173958258667
$ws.Cells.Item(11, 18).Value = 40.258565624328
$ws.Cells.Item(11, 19).Value = 0.0006089676850422485
$ws.Cells.Item(11, 20).Value = 0.0006373492274495819

# Row 12
$ws.Cells.Item(12, 1).Value = 'M2'
$ws.Cells.Item(12, 2).Value = 'Cfh'
$ws.Cells.Item(12, 3).Value = 'Sell'
$ws.Cells.Item(12, 4).Value = 'M1'
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 30.48803466666666
$ws.Cells.Item(12, 8).Value = 91.46410399999999
$ws.Cells.Item(12, 9).Value = 0.2751329786452017
$ws.Cells.Item(12, 10).Value = 0.2879558237532028
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 28.93198366666667
$ws.Cells.Item(12, 14).Value = 86.795951
$ws.Cells.Item(12, 15).Value = 0.4364589970461662
$ws.Cells.Item(12, 16).Value = 0.4364589970461662
$ws.Cells.Item(12, 17).Value = 882.079321004767
$ws.Cells.Item(12, 18).Value = 7938.713889042903
$ws.Cells.Item(12, 19).Value = 0.120084263913809
$ws.Cells.Item(12, 20).Value = 0.1256809100289255

# Row 13
$ws.Cells.Item(13, 1).Value = 'M2'
$ws.Cells.Item(13, 2).Value = 'Cfh'
$ws.Cells.Item(13, 3).Value = 'Sell'
$ws.Cells.Item(13, 4).Value = 'M2'
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 30.48803466666666
$ws.Cells.Item(13, 8).Value = 91.46410399999999
$ws.Cells.Item(13, 9).Value = 0.2751329786452017
$ws.Cells.Item(13, 10).Value = 0.2879558237532028
$ws.Cells.Item(13, 11).Value = 2.0
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 37.20927366666667
$ws.Cells.Item(13, 14).Value = 111.627821
$ws.Cells.Item(13, 15).Value = 0.5613276452965987
$ws.Cells.Item(13, 16).Value = 0.5613276452965988
$ws.Cells.Item(13, 17).Value = 1134.43762547082
$ws.Cells.Item(13, 18).Value = 10209.93862923738
$ws.Cells.Item(13, 19).Value = 0.1544397470463504
$ws.Cells.Item(13, 20).Value = 0.1616375644968278

# Row 14
$ws.Cells.Item(14, 1).Value = 'sCs'
$ws.Cells.Item(14, 2).Value = 'Cfh'
$ws.Cells.Item(14, 3).Value = 'Sell'
$ws.Cells.Item(14, 4).Value = 'ECs'
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 14.8035735
$ws.Cells.Item(14, 8).Value = 29.607147
$ws.Cells.Item(14, 9).Value = 0.1335917948197964
$ws.Cells.Item(14, 10).Value = 0.09321198186522625
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 0.146719
$ws.Cells.Item(14, 14).Value = 0.440157
$ws.Cells.Item(14, 15).Value = 0.002213357657235064
$ws.Cells.Item(14, 16).Value = 0.002213357657235064
$ws.Cells.Item(14, 17).Value = 2.1719655003465
$ws.Cells.Item(14, 18).Value = 13.031793002079
$ws.Cells.Item(14, 19).Value = 0.000295686422008172
$ws.Cells.Item(14, 20).Value = 0.0002063114538074544

# Row 15
$ws.Cells.Item(15, 1).Value = 'sCs'
$ws.Cells.Item(15, 2).Value = 'Cfh'
$ws.Cells.Item(15, 3).Value = 'Sell'
$ws.Cells.Item(15, 4).Value = 'M1'
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 14.8035735
$ws.Cells.Item(15, 8).Value = 29.607147
$ws.Cells.Item(15, 9).Value = 0.1335917948197964
$ws.Cells.Item(15, 10).Value = 0.09321198186522625
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 28.93198366666667
$ws.Cells.Item(15, 14).Value = 86.795951
$ws.Cells.Item(15, 15).Value = 0.4364589970461662
$ws.Cells.Item(15, 16).Value = 0.4364589970461662
$ws.Cells.Item(15, 17).Value = 428.2967467102995
$ws.Cells.Item(15, 18).Value = 2569.780480261797
$ws.Cells.Item(15, 19).Value = 0.05830734078064558
$ws.Cells.Item(15, 20).Value = 0.04068320811758208

# Row 16
$ws.Cells.Item(16, 1).Value = 'sCs'
$ws.Cells.Item(16, 2).Value = 'Cfh'
$ws.Cells.Item(16, 3).Value = 'Sell'
$ws.Cells.Item(16, 4).Value = 'M2'
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 14.8035735
$ws.Cells.Item(16, 8).Value = 29.607147
$ws.Cells.Item(16, 9).Value = 0.1335917948197964
$ws.Cells.Item(16, 10).Value = 0.09321198186522625
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 37.20927366666667
$ws.Cells.Item(16, 14).Value = 111.627821
$ws.Cells.Item(16, 15).Value = 0.5613276452965987
$ws.Cells.Item(16, 16).Value = 0.5613276452965988
$ws.Cells.Item(16, 17).Value = 550.8302176061145
$ws.Cells.Item(16, 18).Value = 3304.981305636687
$ws.Cells.Item(16, 19).Value = 0.0749887676171427
$ws.Cells.Item(16, 20).Value = 0.05232246229383673

